$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:E1) already contains: Full Name, Phone Number, Medication, Quantity, Shipping Address.
# Append the order rows submitted via the intake form (rows 2-5).
# Phone Number column is formatted as text so values like "3234" stay literal
# instead of being coerced to numbers.

$ws.Range("B2:B5").NumberFormat = "@"

$ws.Range("A2").Value = "Chirayu Sahu"
$ws.Range("B2").Value = "3234"
$ws.Range("C2").Value = "paracetamol"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Vit Vellore"

$ws.Range("A3").Value = "Chirayu Sahu"
$ws.Range("B3").Value = "8770195578"
$ws.Range("C3").Value = "paracetamol"
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = "Vit Vellore"

$ws.Range("A4").Value = "Chirayu Sahu"
$ws.Range("B4").Value = "8770195578"
$ws.Range("C4").Value = "paracetamol"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Vit Vellore"

$ws.Range("A5").Value = "Chirayu Sahu"
$ws.Range("B5").Value = "8770195578"
$ws.Range("C5").Value = "aspirin"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "Vit Vellore"
